$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LOM3223")

# Helper: write plain text into a cell without Excel's "smart" numeric/date
# auto-conversion kicking in (e.g. "01/01/2023" would otherwise be stored as
# a date serial). Going through a formula -> copy -> paste-values round trip
# keeps the cell's existing style/number-format untouched.
function Set-PlainText($range, [string]$text) {
    $escaped = $text -replace '"', '""'
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)  # xlPasteValues
    $excel.CutCopyMode = $false
}

# Helper: copy just the formatting (style) from one cell/range to another,
# the same way dragging the format painter or Paste-Special-Formats would.
function Copy-CellFormat($sourceRange, $destRange) {
    $sourceRange.Copy()
    $destRange.PasteSpecial(-4122)  # xlPasteFormats
    $excel.CutCopyMode = $false
}

# 1. "Ativação:" date text updated from 01/01/2012 to 01/01/2023.
#    The same shared string is reused verbatim by the "Programa resumido:"
#    row (B13/C13), so both places change together.
Set-PlainText $ws.Range("B8")  "01/01/2023"
Set-PlainText $ws.Range("C8")  "01/01/2023"
Set-PlainText $ws.Range("B13") "01/01/2023"
Set-PlainText $ws.Range("C13") "01/01/2023"

# 2. "Objectives:" row (row 11) gains an English translation in B/C.
Copy-CellFormat $ws.Range("B10") $ws.Range("B11")
Copy-CellFormat $ws.Range("C10") $ws.Range("C11")
$ws.Range("B11").Value = "Provide the student with the basic knowledge of magnetic and superconducting materials aiming their application in devices."
$ws.Range("C11").Value = "Provide the student with the basic knowledge of magnetic and superconducting materials aiming their application in devices."

# 3. "Short syllabus:" row (row 14) and "Syllabus:" row (row 16) gain the
#    same English text in B/C.
$syllabusText = "Fundamental concepts of magnetic properties of matter. Electron magnetism. Ferromagnetism. Magnetic Materials and Applications: soft and hard. Exchange interaction in oxides and metals. Magnetism - Classical Phenomenology: diamagnetism and paramagnetism. Magnetism - Quantum Phenomenology: ferromagnetism. Magnetic Anisotropy and Spin-Orbit Interaction. Magnetostriction and magnetostrictive materials -Introduction and applications. Basic concepts of superconductivity. Superconductivity - Quantum Origin. Superwave – Consequences. Quantum Interference – SQUID. Superconducting Materials and Applications"

Copy-CellFormat $ws.Range("B15") $ws.Range("B14")
Copy-CellFormat $ws.Range("C15") $ws.Range("C14")
$ws.Range("B14").Value = $syllabusText
$ws.Range("C14").Value = $syllabusText

Copy-CellFormat $ws.Range("B15") $ws.Range("B16")
Copy-CellFormat $ws.Range("C15") $ws.Range("C16")
$ws.Range("B16").Value = $syllabusText
$ws.Range("C16").Value = $syllabusText

# 4. "Norma de recuperação:" row (row 20) text replaced.
$ws.Range("B20").Value = "A nota final , antes da recuperação é dada pela média aritmética das notas das avaliações escritas e da nota do seminário apresentado, se aplicável."
$ws.Range("C20").Value = "A nota final , antes da recuperação é dada pela média aritmética das notas das avaliações escritas e da nota do seminário apresentado, se aplicável."
